$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.042.36"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "2.332.41"
$ws.Range("E3").Value = "  +1.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.82%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.47%  "

# Row 12
$ws.Range("E12").Value = "  -0.29%  "

# Row 13
$ws.Range("E13").Value = "  +1.65%  "

# Row 14
$ws.Range("E14").Value = "  -1.96%  "

# Row 15
$ws.Range("D15").Value = "2.693.02"
$ws.Range("E15").Value = "  +1.09%  "

# Row 16
$ws.Range("D16").Value = "2.332.49"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17
$ws.Range("E17").Value = "  +1.33%  "

# Row 18
$ws.Range("D18").Value = "42.978.28"
$ws.Range("E18").Value = "  +0.19%  "

# Row 19
$ws.Range("E19").Value = "  -2.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.59%  "

# Row 21
$ws.Range("E21").Value = "  -0.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("E24").Value = "  +4.25%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.00%  "

# Row 28
$ws.Range("E28").Value = "  -5.80%  "

# Row 29
$ws.Range("E29").Value = "  +1.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.75%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "142.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.21%  "

# Row 32
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
$ws.Range("E33").Value = "  +0.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "

# Row 35
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0703"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.46%  "

# Row 38
$ws.Range("E38").Value = "  -2.43%  "

# Row 39
$ws.Range("E39").Value = "  +0.28%  "

# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.47%  "

# Row 41
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "

# Row 42
$ws.Range("E42").Value = "  -0.37%  "

# Row 43
$ws.Range("D43").Value = "1.935.07"
$ws.Range("E43").Value = "  -3.19%  "

# Row 44
$ws.Range("E44").Value = "  -0.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.85%  "

# Row 46
$ws.Range("E46").Value = "  -2.75%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.36%  "

# Row 48
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.561.09"
$ws.Range("E48").Value = "  +1.15%  "

# Row 49
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.80%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.77%  "
